$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new cell R3, same style as Q3 (s="9"), stays empty ---
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null

# --- Row 4: new cell R4 = 2021, same style as Q4 (s="16") ---
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 2021

# --- Row 5: R5 = 0.9, new style (font 3 / no border, vertical-top alignment) ---
$ws.Range("P5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").VerticalAlignment = -4160
$ws.Range("R5").Value = 0.9

# --- Row 6: R6 = 6.5, same style as Q6 (s="17") ---
$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 6.5

$excel.CutCopyMode = 0

# --- Move the active selection to T5, matching the new sheetView selection ---
$ws.Range("T5").Select() | Out-Null
